$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 616.2
$ws.Range("I38").Value = 18
$ws.Range("J38").Value = 3009
$ws.Range("K38").Value = 54
$ws.Range("L38").Value = 9027
$ws.Range("M38").Value = 318
$ws.Range("N38").Value = -9771

$ws.Range("H53").Value = 665.1818
$ws.Range("I53").Value = 595.5714
$ws.Range("K53").Value = 595.5714
$ws.Range("M53").Value = 41.42859999999996

$ws.Range("H113").Value = 10666.333
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -6746

$ws.Range("H138").Value = 5428.7144
$ws.Range("I138").Value = 6274.8
$ws.Range("J138").Value = 3313.5
$ws.Range("K138").Value = 18824.4
$ws.Range("L138").Value = 9940.5
$ws.Range("M138").Value = -13684.4
$ws.Range("N138").Value = -20220.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1174.25
$ws.Range("I2").Value = 1232.3334
$ws.Range("K2").Value = 1232.3334
$ws.Range("M2").Value = -1119.3334

$ws.Range("H53").Value = 10038.5
$ws.Range("I53").Value = 10038.5
$ws.Range("K53").Value = 10038.5
$ws.Range("M53").Value = -9356.5

$ws.Range("H55").Value = 49472.75
$ws.Range("J55").Value = 49472.75
$ws.Range("L55").Value = 49472.75
$ws.Range("N55").Value = -50102.75

$ws.Range("H74").Value = 2085.3
$ws.Range("I74").Value = 2108.4736
$ws.Range("K74").Value = 2108.4736
$ws.Range("M74").Value = -1234.4736

$ws.Range("H77").Value = 2085.3
$ws.Range("I77").Value = 2108.4736
$ws.Range("K77").Value = 10542.368
$ws.Range("M77").Value = -6174.367999999999

$ws.Range("H97").Value = 1778.0834
$ws.Range("I97").Value = 1159.7778
$ws.Range("J97").Value = 3633
$ws.Range("K97").Value = 1159.7778
$ws.Range("L97").Value = 3633
$ws.Range("M97").Value = -663.7778000000001
$ws.Range("N97").Value = -4625

$ws.Range("H101").Value = 29150.5
$ws.Range("J101").Value = 29150.5
$ws.Range("L101").Value = 29150.5
$ws.Range("N101").Value = -35640.5

$ws.Range("H116").Value = 1174.25
$ws.Range("I116").Value = 1232.3334
$ws.Range("K116").Value = 1232.3334
$ws.Range("M116").Value = 1061.6666

$ws.Range("H122").Value = 5753.5
$ws.Range("J122").Value = 3007
$ws.Range("L122").Value = 9021
$ws.Range("N122").Value = -13921

$ws.Range("H132").Value = 144810.28
$ws.Range("I132").Value = 201661.4
$ws.Range("K132").Value = 604984.2
$ws.Range("M132").Value = -602454.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1174.25
$ws.Range("I3").Value = 1232.3334
$ws.Range("K3").Value = 1232.3334
$ws.Range("M3").Value = -1118.3334

$ws.Range("H86").Value = 1401.25
$ws.Range("I86").Value = 368.66666
$ws.Range("J86").Value = 4499
$ws.Range("K86").Value = 368.66666
$ws.Range("L86").Value = 4499
$ws.Range("M86").Value = 754.33334
$ws.Range("N86").Value = -6745

$ws.Range("H89").Value = 1401.25
$ws.Range("I89").Value = 368.66666
$ws.Range("J89").Value = 4499
$ws.Range("K89").Value = 1843.3333
$ws.Range("L89").Value = 22495
$ws.Range("M89").Value = 3772.6667
$ws.Range("N89").Value = -33727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1875
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1150

$ws.Range("H32").Value = 2585
$ws.Range("I32").Value = 2585
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2585
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2269
$ws.Range("N32").ClearContents()

$ws.Range("H58").Value = 1068.0834
$ws.Range("I58").Value = 1161.7
$ws.Range("K58").Value = 1161.7
$ws.Range("M58").Value = -958.7

$ws.Range("H122").Value = 2580.4167
$ws.Range("I122").Value = 511.125
$ws.Range("K122").Value = 1533.375
$ws.Range("M122").Value = 916.625

$ws.Range("H136").Value = 1068.0834
$ws.Range("I136").Value = 1161.7
$ws.Range("K136").Value = 3485.1
$ws.Range("M136").Value = -935.1000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 111
$ws.Range("I38").Value = 137.5
$ws.Range("J38").Value = 93.333336
$ws.Range("K38").Value = 412.5
$ws.Range("L38").Value = 280.000008
$ws.Range("M38").Value = -65.5
$ws.Range("N38").Value = -974.000008

$ws.Range("H57").Value = 22499.5
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 60000
$ws.Range("N57").Value = -61118

$ws.Range("H103").Value = 175.8
$ws.Range("I103").Value = 25
$ws.Range("J103").Value = 213.5
$ws.Range("K103").Value = 75
$ws.Range("L103").Value = 640.5
$ws.Range("M103").Value = 804
$ws.Range("N103").Value = -2398.5

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H139").Value = 56498.89
$ws.Range("I139").Value = 2122.5
$ws.Range("J139").Value = 100000
$ws.Range("K139").Value = 6367.5
$ws.Range("L139").Value = 300000
$ws.Range("M139").Value = -1227.5
$ws.Range("N139").Value = -310280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 113466.664
$ws.Range("I3").Value = 500075
$ws.Range("K3").Value = 500075
$ws.Range("M3").Value = -499959

$ws.Range("H26").Value = 24980
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 24980
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 24980
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -25540

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H49").Value = 15030
$ws.Range("J49").Value = 15030
$ws.Range("L49").Value = 15030
$ws.Range("N49").Value = -15398

$ws.Range("H50").Value = 24980
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 24980
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 24980
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -25976

$ws.Range("H97").Value = 456.25
$ws.Range("J97").Value = 175
$ws.Range("L97").Value = 175
$ws.Range("N97").Value = -1167

$ws.Range("H122").Value = 4368.294
$ws.Range("I122").Value = 3208.7273
$ws.Range("K122").Value = 9626.1819
$ws.Range("M122").Value = -7176.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1220.4
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -505

$ws.Range("H27").Value = 1220.4
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 800
$ws.Range("M27").Value = -693

$ws.Range("H40").Value = 1691666.6
$ws.Range("I40").Value = 35000
$ws.Range("K40").Value = 35000
$ws.Range("M40").Value = -34864

$ws.Range("H46").Value = 999.2
$ws.Range("I46").Value = 999.2
$ws.Range("K46").Value = 999.2
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -811.2

$ws.Range("H55").Value = 757.625
$ws.Range("I55").Value = 940.75
$ws.Range("J55").Value = 574.5
$ws.Range("K55").Value = 940.75
$ws.Range("L55").Value = 574.5
$ws.Range("M55").Value = -767.75
$ws.Range("N55").Value = -920.5

$ws.Range("H68").Value = 4656.857
$ws.Range("I68").Value = 3999.6667
$ws.Range("J68").Value = 5149.75
$ws.Range("K68").Value = 3999.6667
$ws.Range("L68").Value = 5149.75
$ws.Range("M68").Value = -3250.6667
$ws.Range("N68").Value = -6647.75

$ws.Range("H71").Value = 4656.857
$ws.Range("I71").Value = 3999.6667
$ws.Range("J71").Value = 5149.75
$ws.Range("K71").Value = 19998.3335
$ws.Range("L71").Value = 25748.75
$ws.Range("M71").Value = -16254.3335
$ws.Range("N71").Value = -33236.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
